$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("E1").Value = "Tag"

# "vip" group (entered together -> single new shared string, index right after Tag)
$ws.Range("E2").Value = "vip"
$ws.Range("E5").Value = "vip"
$ws.Range("E8").Value = "vip"
$ws.Range("E11").Value = "vip"
$ws.Range("E15").Value = "vip"
$ws.Range("E19").Value = "vip"

# "old" group
$ws.Range("E6").Value = "old"
$ws.Range("E9").Value = "old"
$ws.Range("E17").Value = "old"
$ws.Range("E21").Value = "old"

# "company" group
$ws.Range("E3").Value = "company"
$ws.Range("E14").Value = "company"

# "new" group
$ws.Range("E10").Value = "new"
$ws.Range("E13").Value = "new"
$ws.Range("E18").Value = "new"
$ws.Range("E20").Value = "new"

$ws.Range("H6").Select() | Out-Null
